# ArrayVsDictionaryBenchmark.xlsx -- refresh the "String Concatention" sheet
# with a re-run of the BenchmarkDotNet results (new PlusOperator vs
# StringMutate/Concat numbers) and add the Gen0 / Allocated columns that
# BenchmarkDotNet now reports, then resize the chart to make room for the
# wider legend.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("String Concatention")

# ---------------------------------------------------------------------
# 1. Refresh the benchmark timings in B2:J8 with the new run's numbers.
# ---------------------------------------------------------------------
$newData = @(
    @(30.31, 43.51, 63.08, 90.85, 103.37, 123.32, 135.19, 148.81, 181.85),
    @(31.97, 68.36, 124.42, 165.77, 238.82, 331.89, 384.8, 464.51, 550.91),
    @(30.14, 42.61, 60.26, 74.87, 92.47, 105.45, 113.45, 125.57, 154.71),
    @(31.92, 68.67, 117.01, 171.78, 228.46, 295.26, 373.92, 457.44, 564.54999999999995),
    @(116.05, 167.63, 203.84, 255.88, 273.02, 284.45999999999998, 303.91000000000003, 386.91, 405.28),
    @(116.51, 166.51, 190.44, 251.62, 287.89, 282.52, 317.37, 394, 405.89),
    @(44.76, 64.45, 75.22, 92.78, 114.59, 120.79, 136.43, 151.74, 182.27)
)

$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")
for ($r = 0; $r -lt $newData.Count; $r++) {
    $row = $r + 2
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range($cols[$c] + $row).Value = $newData[$r][$c]
    }
}

# ---------------------------------------------------------------------
# 2. Add the new "Gen0" / "Alloc" columns (K, L) that BenchmarkDotNet now
#    reports alongside the mean timings.
# ---------------------------------------------------------------------
$ws.Range("K1").Value = "Gen0"
$ws.Range("L1").Value = "Alloc"
$ws.Range("B1:J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null

$genAlloc = @(
    @(0.27060000000000001, "1136 B"),
    @(1.3513999999999999, "5672 B"),
    @(0.24579999999999999, "1032 B"),
    @(1.3513999999999999, "5672 B"),
    @(0.7319, "3072 B"),
    @(0.7319, "3072 B"),
    @(0.24959999999999999, "1048 B")
)

for ($r = 0; $r -lt $genAlloc.Count; $r++) {
    $row = $r + 2
    $ws.Range("B" + $row).Copy()
    $ws.Range("K" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("K" + $row).Value = $genAlloc[$r][0]

    $ws.Range("L" + $row).Value = $genAlloc[$r][1]
    $ws.Range("B" + $row).Copy()
    $ws.Range("L" + $row).PasteSpecial(-4122) | Out-Null
    $ws.Range("L" + $row).Value = $genAlloc[$r][1]
    $ws.Range("L" + $row).HorizontalAlignment = -4152
}

$ws.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Resize / reposition the line chart so the wider legend (it now has
#    to make room for the extra series labels) still fits nicely.
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$co.Left = 754.5537109375
$co.Top = 2.6251181102362207
$co.Width = 861.625
$co.Height = 571.8749606299212

$legend = $co.Chart.Legend
$legend.Left = 0.13951558938424077
$legend.Top = 0.1118787167340026
$legend.Width = 0.46152401638786217
$legend.Height = 0.31895746884143572

# ---------------------------------------------------------------------
# 4. Print setup + make "String Concatention" the active sheet/tab, with
#    the last selection left on AC17 (mirrors the author's saved state).
# ---------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

$ws.Activate()
$ws.Range("AC17").Select() | Out-Null
